$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
Write-Host "Window type:" $win.GetType().FullName
Write-Host ("ScrollColumn before: " + $win.ScrollColumn)
Write-Host ("ScrollRow before: " + $win.ScrollRow)
$win.ScrollColumn = 3
$win.ScrollRow = 4
Write-Host ("ScrollColumn after: " + $win.ScrollColumn)
Write-Host ("ScrollRow after: " + $win.ScrollRow)
